$d = $word.ActiveDocument

$p = $d.Paragraphs(3).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Mux2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>nBit</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(8).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:r><w:t>el</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve"> bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(14).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Entity name:</w:t></w:r><w:r><w:t xml:space="preserve"> Reg</w:t></w:r><w:r><w:t xml:space="preserve"> =&gt; Generic n := 32</w:t></w:r></w:p>')

$p = $d.Paragraphs(16).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Rst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(17).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Clk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(19).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>RegInput</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/><w:t>(</w:t></w:r><w:r><w:t>n-1</w:t></w:r><w:r><w:t xml:space="preserve"> bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(21).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>RegOutput</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/><w:t>(32 bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(25).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Entity name:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RegFile</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(27).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Rst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(28).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Clk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(29).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>WriteAddress</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/><w:t>(3 bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(32).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>WriteData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(32 bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(33).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>readEnable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(34).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>writeEnable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(51).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>opCo</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t>e</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(3 bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(55).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Cout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1 bit)</w:t></w:r></w:p>')

$p = $d.Paragraphs(58).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Entity name:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>flagControl</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p = $d.Paragraphs(60).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>alu</w:t></w:r><w:r><w:t>R</w:t></w:r><w:r><w:t>es</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(3</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> bits)</w:t></w:r></w:p>')

$p = $d.Paragraphs(62).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Setc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve"> bits)</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>#control signal</w:t></w:r><w:r><w:t xml:space="preserve"> to set carry</w:t></w:r></w:p>')

$p = $d.Paragraphs(98).Range
$r = $d.Range($p.Start, $p.End)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">//should be address in buffer after decode since if second instruction is waiting for data in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>alu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, it has to wait tell it reaches buffer after </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>alu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and in this case the second instruction will have reached buffer after decode</w:t></w:r></w:p>')

Write-Output "DONE"